$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 274 (this shifts the existing rows 274:296
# down to 275:297, preserving all their data/formatting, and carries the
# row-274 date-column style down onto the new blank row).
$ws.Rows("274:274").Insert()

# Populate the newly inserted row 274 with the new weekly price record.
$ws.Cells.Item(274, 1).Value = 11
$ws.Cells.Item(274, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(274, 3).Value = "Bíobío"
$ws.Cells.Item(274, 4).Value = 45013
$ws.Cells.Item(274, 5).Value = 8
$ws.Cells.Item(274, 6).Value = 100112040
$ws.Cells.Item(274, 7).Value = "Cilantro"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 60
$ws.Cells.Item(274, 11).Value = 6500
$ws.Cells.Item(274, 12).Value = 7000
$ws.Cells.Item(274, 13).Value = 6750
$ws.Cells.Item(274, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(274, 15).Value = "Región Metropolitana"
$ws.Cells.Item(274, 16).Value = 188
$ws.Cells.Item(274, 17).Value = 36
$ws.Cells.Item(274, 18).Value = "Hortaliza"
